$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.01470429199441242
$ws.Range("B3").Value = -0.04558021581030089
$ws.Range("B4").Value = -0.009826403510106795
$ws.Range("B5").Value = 0.008226121485872784
$ws.Range("B6").Value = -0.04316008579431754
$ws.Range("B7").Value = -0.01157253281283423
$ws.Range("B8").Value = -0.006493037314404795
$ws.Range("B9").Value = -0.02572169932097143
$ws.Range("B10").Value = -0.08015707014535689
$ws.Range("B11").Value = 0.06473596093476565
$ws.Range("B12").Value = 0.07139121916325228
$ws.Range("B13").Value = 0.07733637174846859
$ws.Range("B14").Value = 0.1233898754576115
$ws.Range("B15").Value = 0.2091587838646004
$ws.Range("B16").Value = -0.08353428163828867
$ws.Range("B17").Value = -0.002495154701341215
$ws.Range("B18").Value = 0.01671036800076502
$ws.Range("B19").Value = 0.2753602545626601
$ws.Range("B20").Value = -0.222677881257742
$ws.Range("B21").Value = -0.01930698667692329
$ws.Range("B22").Value = -0.04105922513277995
$ws.Range("B23").Value = -0.141938455995426
$ws.Range("B24").Value = 0.1048743112924959
$ws.Range("B25").Value = 0.04424339686558416
$ws.Range("B26").Value = 0.05949019800012678
$ws.Range("B27").Value = 0.01861251411408677
$ws.Range("B28").Value = 0.07736807577675232
$ws.Range("B29").Value = -0.0250016481418028
$ws.Range("B30").Value = -0.01961362217271759
$ws.Range("B31").Value = -0.0294715145912881
$ws.Range("B32").Value = -0.09576921519656952
$ws.Range("B33").Value = 0.006417395024311053
$ws.Range("B34").Value = -0.02292509128313051
$ws.Range("B35").Value = 0.06332293387393957
$ws.Range("B36").Value = 0.03937713790365305
$ws.Range("B37").Value = 0.07042443893857833
$ws.Range("B38").Value = 0.06768668814173028
$ws.Range("B39").Value = 0.02909378411594531
$ws.Range("B40").Value = 0.1119119236831914
$ws.Range("B41").Value = 0.1160916927122154
$ws.Range("B42").Value = 0.079998289855468
$ws.Range("B43").Value = 0.0295060767461054
$ws.Range("B44").Value = 0.09388609445062236
$ws.Range("B45").Value = 0.02593231010365504
$ws.Range("B46").Value = 0.05303268938139529
$ws.Range("B47").Value = -0.03745826781061814
$ws.Range("B48").Value = 0.06270911510394725
$ws.Range("B49").Value = 0.04609120822202557
$ws.Range("B50").Value = -0.05310395082503255
$ws.Range("B51").Value = -0.01480859739746029
